# Update crypto price/volume data in the worksheet.
# D-column "Price" values are plain text strings in the source data; some of them
# are simple decimal numbers (e.g. "51.13") which Excel would otherwise silently
# coerce to a numeric cell type on assignment. Force those to stay text by pre-
# setting the cell to the Text number format, mirroring how Excel treats a value
# typed with a leading apostrophe.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Range, $Text) {
    if ($Text -match '^-?[0-9]+(\.[0-9]+)?$') {
        $Range.NumberFormat = "@"
    }
    $Range.Value = $Text
}

Set-TextValue $ws.Range("D2") "23.697.69"
$ws.Range("E2").Value = "  +0.82%  "
Set-TextValue $ws.Range("D3") "1.655.44"
$ws.Range("E3").Value = "  +0.74%  "
Set-TextValue $ws.Range("D4") "1.002"
$ws.Range("E4").Value = "  +0.22%  "
Set-TextValue $ws.Range("D5") "1.001"
$ws.Range("E5").Value = "  +0.10%  "
Set-TextValue $ws.Range("D6") "302.62"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("E7").Value = "  +0.33%  "
Set-TextValue $ws.Range("D8") "0.3607"
$ws.Range("E8").Value = "  -0.29%  "
Set-TextValue $ws.Range("D9") "51.13"
$ws.Range("E9").Value = "  -1.69%  "
Set-TextValue $ws.Range("D10") "0.08190"
$ws.Range("E10").Value = "  -0.77%  "
Set-TextValue $ws.Range("D11") "1.229"
$ws.Range("E11").Value = "  -0.44%  "
Set-TextValue $ws.Range("D12") "1.002"
$ws.Range("E12").Value = "  +0.14%  "
Set-TextValue $ws.Range("D13") "22.43"
$ws.Range("E13").Value = "  -0.66%  "
Set-TextValue $ws.Range("D14") "6.443"
$ws.Range("E14").Value = "  -0.44%  "
Set-TextValue $ws.Range("D15") "7.431"
$ws.Range("E15").Value = "  +0.97%  "
Set-TextValue $ws.Range("D16") "0.00001223"
$ws.Range("E16").Value = "  -1.49%  "
Set-TextValue $ws.Range("D17") "1.654.99"
$ws.Range("E17").Value = "  +0.94%  "
Set-TextValue $ws.Range("D18") "97.35"
$ws.Range("E18").Value = "  +2.29%  "
Set-TextValue $ws.Range("D19") "0.07029"
$ws.Range("E19").Value = "  +0.88%  "
Set-TextValue $ws.Range("D20") "6.790"
$ws.Range("E20").Value = "  +2.77%  "
Set-TextValue $ws.Range("D21") "17.56"
$ws.Range("E21").Value = "  +0.17%  "
Set-TextValue $ws.Range("D22") "1.002"
$ws.Range("E22").Value = "  +0.16%  "
Set-TextValue $ws.Range("D23") "12.73"
$ws.Range("E23").Value = "  +1.50%  "
Set-TextValue $ws.Range("D24") "23.710.13"
$ws.Range("E24").Value = "  +0.92%  "
Set-TextValue $ws.Range("D25") "2.503"
$ws.Range("E25").Value = "  -1.19%  "
Set-TextValue $ws.Range("D26") "3.020"
$ws.Range("E26").Value = "  -2.04%  "
Set-TextValue $ws.Range("D27") "21.24"
$ws.Range("E27").Value = "  +0.18%  "
Set-TextValue $ws.Range("D28") "154.00"
$ws.Range("E28").Value = "  +1.48%  "
Set-TextValue $ws.Range("D29") "5.233"
$ws.Range("E29").Value = "  -0.75%  "
Set-TextValue $ws.Range("D30") "134.08"
$ws.Range("E30").Value = "  +0.51%  "
Set-TextValue $ws.Range("D31") "1.842.23"
$ws.Range("E31").Value = "  +1.23%  "
Set-TextValue $ws.Range("D32") "7.169"
$ws.Range("E32").Value = "  +9.18%  "
Set-TextValue $ws.Range("D33") "2.247"
$ws.Range("E33").Value = "  +4.06%  "
Set-TextValue $ws.Range("D34") "11.98"
$ws.Range("E34").Value = "  +4.23%  "
Set-TextValue $ws.Range("D35") "1.054"
$ws.Range("E35").Value = "  -3.38%  "
Set-TextValue $ws.Range("D36") "0.02800"
$ws.Range("E36").Value = "  +0.86%  "
Set-TextValue $ws.Range("D37") "0.2503"
$ws.Range("E37").Value = "  -0.59%  "
$ws.Range("E38").Value = "  +0.23%  "
Set-TextValue $ws.Range("D39") "6.094"
$ws.Range("E39").Value = "  +1.74%  "
Set-TextValue $ws.Range("D40") "0.06991"
$ws.Range("E40").Value = "  -0.74%  "
Set-TextValue $ws.Range("D41") "12.98"
$ws.Range("E41").Value = "  +5.51%  "
Set-TextValue $ws.Range("D42") "0.6977"
$ws.Range("E42").Value = "  -1.22%  "
Set-TextValue $ws.Range("D43") "1.327"
$ws.Range("E43").Value = "  -1.61%  "
Set-TextValue $ws.Range("D44") "16.04"
$ws.Range("E44").Value = "  +2.51%  "
Set-TextValue $ws.Range("D45") "0.6504"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("E46").Value = "  +0.13%  "
Set-TextValue $ws.Range("D47") "2.298"
$ws.Range("E47").Value = "  -0.11%  "
Set-TextValue $ws.Range("D48") "3.959"
$ws.Range("E48").Value = "  -0.14%  "
Set-TextValue $ws.Range("D49") "0.07901"
$ws.Range("E49").Value = "  -1.00%  "
Set-TextValue $ws.Range("D50") "127.99"
$ws.Range("E50").Value = "  -1.03%  "
Set-TextValue $ws.Range("D51") "1.177"
$ws.Range("E51").Value = "  -1.44%  "

Write-Output "Updated $(97) cells."
